$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "51.570.28"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +1.66%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.986.94"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +2.30%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "382.07"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +2.11%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "104.20"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +4.65%  "

$ws.Range("E7").Value = "  +2.24%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +2.46%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "36.83"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +2.91%  "

$ws.Range("E11").Value = "  -0.39%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.0861"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +2.55%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "3.463.51"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +2.55%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "18.51"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +3.26%  "

$ws.Range("E15").Value = "  +4.03%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "2.997.89"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +2.94%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "11.28"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +0.72%  "

$ws.Range("E18").Value = "  +1.43%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "51.614.69"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +1.92%  "

$ws.Range("E20").Value = "  +1.43%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "12.57"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +2.36%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0965"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +1.67%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "70.38"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +2.55%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "267.65"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +1.32%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "3.24"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +3.50%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "8.01"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +0.84%  "

$ws.Range("E27").Value = "  +5.01%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "7.21"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -1.56%  "

$ws.Range("E29").Value = "  -0.07%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "26.14"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +2.98%  "

$ws.Range("E31").Value = "  +0.98%  "

$ws.Range("E32").Value = "  +4.64%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "34.65"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +5.42%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "51.43"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +1.43%  "

$ws.Range("E35").Value = "  +0.53%  "

$ws.Range("E36").Value = "  +2.20%  "

$ws.Range("E37").Value = "  +0.00%  "

$ws.Range("E38").Value = "  +8.15%  "

$ws.Range("E39").Value = "  +3.84%  "

$ws.Range("E40").Value = "  +5.79%  "

$ws.Range("E41").Value = "  +1.70%  "

$ws.Range("E42").Value = "  +3.14%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "3.84"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +14.83%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "122.81"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +3.05%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "21.49"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +3.01%  "

$ws.Range("E46").Value = "  +0.15%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.273"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +2.58%  "

$ws.Range("E48").Value = "  +0.67%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "2.035.63"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +2.89%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "3.283.31"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +2.37%  "

$ws.Range("E51").Value = "  +2.29%  "

